$d = $word.ActiveDocument

function Split-ReplaceTail {
    param(
        [string]$FindText,
        [string]$KeepPrefix,
        [string]$NewTail
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $FindText"
    }

    $tailStart = $rng.Start + $KeepPrefix.Length
    $tailEnd = $rng.End

    # Remove the old tail text, leaving a collapsed range right after the prefix.
    $tailRange = $d.Range($tailStart, $tailEnd)
    $tailRange.Text = ""

    # Insert the new tail as a brand-new run (keeps it split from the prefix run).
    $collapsed = $d.Range($tailStart, $tailStart)
    $collapsed.InsertAfter($NewTail)
}

Split-ReplaceTail "- Deployed via Railway" "- Deployed via " "Render"

Split-ReplaceTail "- Collection URL: [YOUR_PUBLIC_GIST_LINK_HERE]" "- Collection URL: " "https://github.com/ranadive-25/splitwise-backend/blob/main/Split%20App%20-%20DevDynamics.postman_collection.json"

Split-ReplaceTail "- The backend is deployed on Railway" "- The backend is deployed on " "Render"
